$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header columns for the season record (Wins/Losses/Ties),
# matching the bold/centered/bordered style used by the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill the season-record values for every player row (2-55)
for ($r = 2; $r -le 55; $r++) {
    $ws.Cells.Item($r, 30).Value = 62
    $ws.Cells.Item($r, 31).Value = 100
    $ws.Cells.Item($r, 32).Value = 0
}
